# Applies the "Player Info" sheet addition + "ODI Batting" MATCH_CODE rename
# described by the commit's diff:
#   - new worksheet "Player Info" inserted before "ODI Batting", with an
#     ID/NAME/BATTING_HAND/BOWL_STYLE header row and one data row for
#     player 5929 (Joshua Da Silva)
#   - "ODI Batting"!D1 header renamed MATCH_CARD_LINK -> MATCH_CODE, and the
#     full scorecard URLs in D2/D3 replaced with the bare match codes

$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet "Player Info" before the existing "ODI Batting" sheet ---
$odi = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($odi)
$playerInfo.Name = "Player Info"

# NOTE: a worksheet reference captured before Add() can resolve by position
# rather than identity once the sheet collection shifts, so re-fetch
# "ODI Batting" by name now that "Player Info" has been inserted ahead of it.
$odi = $wb.Worksheets.Item("ODI Batting")

# --- 2. Header row: reuse ODI Batting's existing header formatting (bold,
#        centered, top-aligned, thin border) by copying its format over ---
$odi.Range("A1:D1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)   # xlPasteFormats

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# --- 3. Data row; leading apostrophe forces text storage (quote-prefix) so the
#        numeric-looking ID stays a string instead of being coerced to a number ---
$playerInfo.Cells.Item(2, 1).Value = "'5929"
$playerInfo.Cells.Item(2, 2).Value = "Joshua Da Silva"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Does Not Bowl | Unknown"

# --- 4. Rename MATCH_CARD_LINK column to MATCH_CODE and replace URL values with bare match codes ---
$odi.Range("D1").Value = "MATCH_CODE"

$odi.Cells.Item(2, 4).Value = "'4443"
$odi.Cells.Item(3, 4).Value = "'4445"
